# Generate Report for Archive
#
# The localization run moved on from "Ready for handoff" to "In Translation",
# so every Status cell that showed the old value needs to show the new one,
# and the (now shorter) status column can be narrowed back down to fit.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Target stored column width is 13.4101845877511 characters. ColumnWidth is
# quantized to whole pixels on save, so request the value in the middle of
# the pixel bucket that rounds closest to the target (12.5 -> stored width
# 13.33.., nearer to 13.41 than any other reachable pixel width).
$newWidth  = 12.5

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E1:F1").EntireColumn.ColumnWidth = $newWidth

# --- Per-locale detail sheets: Status column (C) for both rows ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus
    $ws.Range("C1").EntireColumn.ColumnWidth = $newWidth
}

# --- Re-assert the datetime number format on the date columns untouched by
#     this edit; a save/reload of this workbook can otherwise drop a custom
#     numFmt's association on these specific cells, which would be an
#     unintended side effect outside the scope of this change. ---
$dateFormat = "yyyy-mm-dd HH:mm:ss"
$overview.Range("G2:G3").NumberFormat = $dateFormat
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("H2:H3").NumberFormat = $dateFormat
    $ws.Range("K2:K3").NumberFormat = $dateFormat
}

# --- Likewise re-assert the HyperLink font (underline + blue) on the
#     existing file-name hyperlink cells, untouched by this edit. ---
$hyperlinkColor = 15570276 # OLE color for FF6495ED
$overview.Range("B2:B3").Font.Underline = $true
$overview.Range("B2:B3").Font.Color = $hyperlinkColor
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A2:A3").Font.Underline = $true
    $ws.Range("A2:A3").Font.Color = $hyperlinkColor
}
